$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Rows 1-3: replace their single values with "0M"
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# 2) Insert 10 new rows right after row 3, each holding one new stat value.
#    We repeatedly insert a fresh row before (what is currently) row 4 and
#    fill it in; inserting the values in reverse order means the final
#    top-to-bottom order matches the desired sequence.
$newValues = @('106', '0.00003', '0.00007', '0.00004', '0.00001', '0.00004', '0.00005', '0.00005', '0.00497', '100.0')
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $t.Rows.Add($t.Rows.Item(4)) | Out-Null
    $t.Cell(4, 1).Range.Text = $newValues[$i]
}

# 3) The last three rows used to hold a whole tab-separated line of stats;
#    collapse each back down to a single value (no tabs).
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.99"
$t.Cell($rowCount - 1, 1).Range.Text = "0"
$t.Cell($rowCount, 1).Range.Text = "98"
